$d = $word.ActiveDocument

$d.Content.Find.Execute("87×82=", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=", 2) | Out-Null
$d.Content.Find.Execute("55×13=", $true, $false, $false, $false, $false, $true, 1, $false, "41×79=", 2) | Out-Null
$d.Content.Find.Execute("76×89=", $true, $false, $false, $false, $false, $true, 1, $false, "97×53=", 2) | Out-Null
$d.Content.Find.Execute("85×13=", $true, $false, $false, $false, $false, $true, 1, $false, "43×80=", 2) | Out-Null
$d.Content.Find.Execute("72×10=", $true, $false, $false, $false, $false, $true, 1, $false, "50×40=", 2) | Out-Null
$d.Content.Find.Execute("13×38=", $true, $false, $false, $false, $false, $true, 1, $false, "60×72=", 2) | Out-Null
$d.Content.Find.Execute("20×99=", $true, $false, $false, $false, $false, $true, 1, $false, "27×21=", 2) | Out-Null
$d.Content.Find.Execute("69×61=", $true, $false, $false, $false, $false, $true, 1, $false, "46×58=", 2) | Out-Null
$d.Content.Find.Execute("43×81=", $true, $false, $false, $false, $false, $true, 1, $false, "53×92=", 2) | Out-Null
$d.Content.Find.Execute("51×29=", $true, $false, $false, $false, $false, $true, 1, $false, "21×35=", 2) | Out-Null
$d.Content.Find.Execute("69×94=", $true, $false, $false, $false, $false, $true, 1, $false, "84×72=", 2) | Out-Null
$d.Content.Find.Execute("91×96=", $true, $false, $false, $false, $false, $true, 1, $false, "90×60=", 2) | Out-Null
$d.Content.Find.Execute("48×36=", $true, $false, $false, $false, $false, $true, 1, $false, "91×61=", 2) | Out-Null
$d.Content.Find.Execute("70×14=", $true, $false, $false, $false, $false, $true, 1, $false, "40×98=", 2) | Out-Null
$d.Content.Find.Execute("33×99=", $true, $false, $false, $false, $false, $true, 1, $false, "65×83=", 2) | Out-Null
$d.Content.Find.Execute("72×93=", $true, $false, $false, $false, $false, $true, 1, $false, "18×93=", 2) | Out-Null
$d.Content.Find.Execute("90×77=", $true, $false, $false, $false, $false, $true, 1, $false, "67×89=", 2) | Out-Null
$d.Content.Find.Execute("57×92=", $true, $false, $false, $false, $false, $true, 1, $false, "74×50=", 2) | Out-Null
$d.Content.Find.Execute("73×19=", $true, $false, $false, $false, $false, $true, 1, $false, "11×59=", 2) | Out-Null
$d.Content.Find.Execute("17×88=", $true, $false, $false, $false, $false, $true, 1, $false, "74×99=", 2) | Out-Null
$d.Content.Find.Execute("88×97=", $true, $false, $false, $false, $false, $true, 1, $false, "10×18=", 2) | Out-Null
$d.Content.Find.Execute("60×33=", $true, $false, $false, $false, $false, $true, 1, $false, "85×61=", 2) | Out-Null
$d.Content.Find.Execute("75×54=", $true, $false, $false, $false, $false, $true, 1, $false, "56×91=", 2) | Out-Null
$d.Content.Find.Execute("65×15=", $true, $false, $false, $false, $false, $true, 1, $false, "66×97=", 2) | Out-Null
$d.Content.Find.Execute("10×13=", $true, $false, $false, $false, $false, $true, 1, $false, "18×43=", 2) | Out-Null
$d.Content.Find.Execute("87×61=", $true, $false, $false, $false, $false, $true, 1, $false, "42×87=", 2) | Out-Null
$d.Content.Find.Execute("54×25=", $true, $false, $false, $false, $false, $true, 1, $false, "57×60=", 2) | Out-Null
$d.Content.Find.Execute("17×18=", $true, $false, $false, $false, $false, $true, 1, $false, "28×96=", 2) | Out-Null
$d.Content.Find.Execute("95×81=", $true, $false, $false, $false, $false, $true, 1, $false, "99×62=", 2) | Out-Null
$d.Content.Find.Execute("95×27=", $true, $false, $false, $false, $false, $true, 1, $false, "16×93=", 2) | Out-Null
$d.Content.Find.Execute("87×57=", $true, $false, $false, $false, $false, $true, 1, $false, "58×84=", 2) | Out-Null
$d.Content.Find.Execute("74×71=", $true, $false, $false, $false, $false, $true, 1, $false, "18×94=", 2) | Out-Null
$d.Content.Find.Execute("30×11=", $true, $false, $false, $false, $false, $true, 1, $false, "71×92=", 2) | Out-Null
$d.Content.Find.Execute("42×94=", $true, $false, $false, $false, $false, $true, 1, $false, "66×93=", 2) | Out-Null
$d.Content.Find.Execute("85×90=", $true, $false, $false, $false, $false, $true, 1, $false, "55×56=", 2) | Out-Null
$d.Content.Find.Execute("24×52=", $true, $false, $false, $false, $false, $true, 1, $false, "79×11=", 2) | Out-Null
$d.Content.Find.Execute("90×74=", $true, $false, $false, $false, $false, $true, 1, $false, "63×84=", 2) | Out-Null
$d.Content.Find.Execute("12×69=", $true, $false, $false, $false, $false, $true, 1, $false, "28×39=", 2) | Out-Null
$d.Content.Find.Execute("54×48=", $true, $false, $false, $false, $false, $true, 1, $false, "74×92=", 2) | Out-Null
$d.Content.Find.Execute("33×54=", $true, $false, $false, $false, $false, $true, 1, $false, "42×98=", 2) | Out-Null
$d.Content.Find.Execute("82×41=", $true, $false, $false, $false, $false, $true, 1, $false, "21×30=", 2) | Out-Null
$d.Content.Find.Execute("90×16=", $true, $false, $false, $false, $false, $true, 1, $false, "74×44=", 2) | Out-Null
$d.Content.Find.Execute("95×59=", $true, $false, $false, $false, $false, $true, 1, $false, "18×67=", 2) | Out-Null
$d.Content.Find.Execute("78×62=", $true, $false, $false, $false, $false, $true, 1, $false, "71×98=", 2) | Out-Null
$d.Content.Find.Execute("73×76=", $true, $false, $false, $false, $false, $true, 1, $false, "62×68=", 2) | Out-Null
$d.Content.Find.Execute("25×50=", $true, $false, $false, $false, $false, $true, 1, $false, "67×72=", 2) | Out-Null
$d.Content.Find.Execute("29×41=", $true, $false, $false, $false, $false, $true, 1, $false, "86×29=", 2) | Out-Null
$d.Content.Find.Execute("57×67=", $true, $false, $false, $false, $false, $true, 1, $false, "77×12=", 2) | Out-Null
$d.Content.Find.Execute("57×72=", $true, $false, $false, $false, $false, $true, 1, $false, "78×40=", 2) | Out-Null
$d.Content.Find.Execute("49×61=", $true, $false, $false, $false, $false, $true, 1, $false, "76×46=", 2) | Out-Null
$d.Content.Find.Execute("79×36=", $true, $false, $false, $false, $false, $true, 1, $false, "13×26=", 2) | Out-Null
$d.Content.Find.Execute("98×72=", $true, $false, $false, $false, $false, $true, 1, $false, "93×65=", 2) | Out-Null
$d.Content.Find.Execute("38×46=", $true, $false, $false, $false, $false, $true, 1, $false, "51×52=", 2) | Out-Null
$d.Content.Find.Execute("14×86=", $true, $false, $false, $false, $false, $true, 1, $false, "35×94=", 2) | Out-Null
$d.Content.Find.Execute("66×17=", $true, $false, $false, $false, $false, $true, 1, $false, "81×16=", 2) | Out-Null
$d.Content.Find.Execute("32×33=", $true, $false, $false, $false, $false, $true, 1, $false, "48×64=", 2) | Out-Null
$d.Content.Find.Execute("99×29=", $true, $false, $false, $false, $false, $true, 1, $false, "39×37=", 2) | Out-Null
$d.Content.Find.Execute("65×17=", $true, $false, $false, $false, $false, $true, 1, $false, "56×95=", 2) | Out-Null
$d.Content.Find.Execute("94×40=", $true, $false, $false, $false, $false, $true, 1, $false, "17×19=", 2) | Out-Null
$d.Content.Find.Execute("49×31=", $true, $false, $false, $false, $false, $true, 1, $false, "25×35=", 2) | Out-Null
$d.Content.Find.Execute("67×12=", $true, $false, $false, $false, $false, $true, 1, $false, "86×42=", 2) | Out-Null
$d.Content.Find.Execute("54×87=", $true, $false, $false, $false, $false, $true, 1, $false, "92×52=", 2) | Out-Null
$d.Content.Find.Execute("29×70=", $true, $false, $false, $false, $false, $true, 1, $false, "22×71=", 2) | Out-Null
$d.Content.Find.Execute("58×46=", $true, $false, $false, $false, $false, $true, 1, $false, "55×60=", 2) | Out-Null
$d.Content.Find.Execute("11×66=", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=", 2) | Out-Null
$d.Content.Find.Execute("47×24=", $true, $false, $false, $false, $false, $true, 1, $false, "23×60=", 2) | Out-Null
$d.Content.Find.Execute("45×86=", $true, $false, $false, $false, $false, $true, 1, $false, "34×86=", 2) | Out-Null
$d.Content.Find.Execute("98×62=", $true, $false, $false, $false, $false, $true, 1, $false, "55×33=", 2) | Out-Null
$d.Content.Find.Execute("40×52=", $true, $false, $false, $false, $false, $true, 1, $false, "74×77=", 2) | Out-Null
$d.Content.Find.Execute("99×42=", $true, $false, $false, $false, $false, $true, 1, $false, "84×91=", 2) | Out-Null
$d.Content.Find.Execute("91×40=", $true, $false, $false, $false, $false, $true, 1, $false, "85×16=", 2) | Out-Null
$d.Content.Find.Execute("91×58=", $true, $false, $false, $false, $false, $true, 1, $false, "89×75=", 2) | Out-Null
$d.Content.Find.Execute("64×58=", $true, $false, $false, $false, $false, $true, 1, $false, "66×77=", 2) | Out-Null
$d.Content.Find.Execute("93×38=", $true, $false, $false, $false, $false, $true, 1, $false, "33×44=", 2) | Out-Null
$d.Content.Find.Execute("89×63=", $true, $false, $false, $false, $false, $true, 1, $false, "25×45=", 2) | Out-Null
$d.Content.Find.Execute("14×77=", $true, $false, $false, $false, $false, $true, 1, $false, "40×38=", 2) | Out-Null
$d.Content.Find.Execute("57×81=", $true, $false, $false, $false, $false, $true, 1, $false, "45×44=", 2) | Out-Null
$d.Content.Find.Execute("43×17=", $true, $false, $false, $false, $false, $true, 1, $false, "62×56=", 2) | Out-Null
$d.Content.Find.Execute("67×96=", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=", 2) | Out-Null
$d.Content.Find.Execute("30×53=", $true, $false, $false, $false, $false, $true, 1, $false, "58×76=", 2) | Out-Null
$d.Content.Find.Execute("95×58=", $true, $false, $false, $false, $false, $true, 1, $false, "53×14=", 2) | Out-Null
$d.Content.Find.Execute("48×48=", $true, $false, $false, $false, $false, $true, 1, $false, "26×27=", 2) | Out-Null
$d.Content.Find.Execute("69×73=", $true, $false, $false, $false, $false, $true, 1, $false, "88×63=", 2) | Out-Null
$d.Content.Find.Execute("34×88=", $true, $false, $false, $false, $false, $true, 1, $false, "76×49=", 2) | Out-Null
$d.Content.Find.Execute("88×53=", $true, $false, $false, $false, $false, $true, 1, $false, "61×58=", 2) | Out-Null
$d.Content.Find.Execute("47×88=", $true, $false, $false, $false, $false, $true, 1, $false, "98×55=", 2) | Out-Null
$d.Content.Find.Execute("59×40=", $true, $false, $false, $false, $false, $true, 1, $false, "77×52=", 2) | Out-Null
$d.Content.Find.Execute("72×68=", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=", 2) | Out-Null
$d.Content.Find.Execute("94×28=", $true, $false, $false, $false, $false, $true, 1, $false, "25×26=", 2) | Out-Null
$d.Content.Find.Execute("10×72=", $true, $false, $false, $false, $false, $true, 1, $false, "46×92=", 2) | Out-Null
$d.Content.Find.Execute("21×75=", $true, $false, $false, $false, $false, $true, 1, $false, "94×52=", 2) | Out-Null
$d.Content.Find.Execute("88×81=", $true, $false, $false, $false, $false, $true, 1, $false, "15×69=", 2) | Out-Null
$d.Content.Find.Execute("53×57=", $true, $false, $false, $false, $false, $true, 1, $false, "51×77=", 2) | Out-Null
$d.Content.Find.Execute("72×95=", $true, $false, $false, $false, $false, $true, 1, $false, "94×42=", 2) | Out-Null
$d.Content.Find.Execute("95×40=", $true, $false, $false, $false, $false, $true, 1, $false, "15×22=", 2) | Out-Null
$d.Content.Find.Execute("76×82=", $true, $false, $false, $false, $false, $true, 1, $false, "20×71=", 2) | Out-Null
$d.Content.Find.Execute("33×64=", $true, $false, $false, $false, $false, $true, 1, $false, "44×70=", 2) | Out-Null
$d.Content.Find.Execute("52×55=", $true, $false, $false, $false, $false, $true, 1, $false, "100×22=", 2) | Out-Null
$d.Content.Find.Execute("62×49=", $true, $false, $false, $false, $false, $true, 1, $false, "46×94=", 2) | Out-Null
$d.Content.Find.Execute("51×67=", $true, $false, $false, $false, $false, $true, 1, $false, "66×86=", 2) | Out-Null
